$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Manager column (D) for rows 2-4 from "Steven Panter" to "Dave Allsop"
$ws.Range("D2:D4").Value = "Dave Allsop"

# Update the selected range/active cell to match the new selection
$ws.Activate()
$ws.Range("D2:D4").Select()
